$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$titleShape = $s.Shapes.Item(1)

# The title currently holds two runs ("Let's " + "Play ...cont.") plus a
# trailing empty end-of-paragraph run. Clear the text frame completely
# (this drops the stale end-paragraph run-properties too) and retype the
# merged sentence as a single run.
$titleShape.TextFrame.DeleteText()
$titleShape.TextFrame.TextRange.Text = "Let’s Play                                        …cont."
